# Applies the "minor revisions to seagrass data entry spreadsheets" edit:
#  1. abundance_data: drop the stale scrolled-in sheetView topLeftCell
#     (re-touching the sheet's window state causes it to be recomputed).
#  2. glossary: the three definition rows under the "abundance_data" and
#     "length_data" section headers get cycled (site_code /
#     location_name / sample_collection_date -> sample_collection_date /
#     site_code / location_name), and the view scrolls down to show the
#     length_data block with A42 selected.

$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet

# --- 1. abundance_data: reset stale window position -----------------------
$abundance = $wb.Worksheets.Item("abundance_data")
$abundance.Activate()

# --- 2. glossary: rotate the three glossary rows in each data-sheet block -
$glossary = $wb.Worksheets.Item("glossary")
$glossary.Activate()

# abundance_data block (rows 32-34)
$glossary.Range("B32").Value = "sample_collection_date"
$glossary.Range("C32").Value = "The date the sample was collected in the field"
$glossary.Range("D32").Value = "date"
$glossary.Range("E32").Value = "YYYY-MM-DD"

$glossary.Range("B33").Value = "site_code"
$glossary.Range("C33").Value = "Enter your 6 character site code. Codes can be found in the standards section of the MarineGEO protocol website: https://marinegeo.github.io/standards/"
$glossary.Range("D33").Value = "text"
$glossary.Range("E33").Value = "XXX-YYY"

$glossary.Range("B34").Value = "location_name"
$glossary.Range("C34").Value = "The name of the specific location where the sample was collected; e.g., Curlew Cay"
$glossary.Range("D34").Value = "text"
$glossary.Range("E34").ClearContents()

# length_data block (rows 41-43)
$glossary.Range("B41").Value = "sample_collection_date"
$glossary.Range("C41").Value = "The date the sample was collected in the field"
$glossary.Range("D41").Value = "date"
$glossary.Range("E41").Value = "YYYY-MM-DD"

$glossary.Range("B42").Value = "site_code"
$glossary.Range("C42").Value = "Enter your 6 character site code. Codes can be found in the standards section of the MarineGEO protocol website: https://marinegeo.github.io/standards/"
$glossary.Range("D42").Value = "text"
$glossary.Range("E42").Value = "XXX-YYY"

$glossary.Range("B43").Value = "location_name"
$glossary.Range("C43").Value = "The name of the specific location where the sample was collected; e.g., Curlew Cay"
$glossary.Range("D43").Value = "text"
$glossary.Range("E43").ClearContents()

# Scroll the glossary view down to the length_data block and land on A42.
$glossary.Range("A42").Select()

# Restore the originally active sheet/tab.
$originalActive.Activate()
